$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as the new row 181 (dated
# 2022-03-22 / serial 44642), pushing all the previously existing rows
# 181-221 down by one (to rows 182-222).
$ws.Rows.Item(181).Insert()

# Copy the boilerplate / descriptive columns (Mercado ID, Mercado, Region,
# Codreg, Categoria ID, Categoria, Variedad, Calidad, Origen, Clasificacion)
# from the row directly below, since they are identical for every record in
# this sheet (same market / product).
$ws.Cells.Item(181,1).Value  = $ws.Cells.Item(182,1).Value()
$ws.Cells.Item(181,2).Value  = $ws.Cells.Item(182,2).Value()
$ws.Cells.Item(181,3).Value  = $ws.Cells.Item(182,3).Value()
$ws.Cells.Item(181,5).Value  = $ws.Cells.Item(182,5).Value()
$ws.Cells.Item(181,6).Value  = $ws.Cells.Item(182,6).Value()
$ws.Cells.Item(181,7).Value  = $ws.Cells.Item(182,7).Value()
$ws.Cells.Item(181,8).Value  = $ws.Cells.Item(182,8).Value()
$ws.Cells.Item(181,9).Value  = $ws.Cells.Item(182,9).Value()
$ws.Cells.Item(181,14).Value = $ws.Cells.Item(182,14).Value()
$ws.Cells.Item(181,15).Value = $ws.Cells.Item(182,15).Value()
$ws.Cells.Item(181,17).Value = $ws.Cells.Item(182,17).Value()
$ws.Cells.Item(181,18).Value = $ws.Cells.Item(182,18).Value()

# New values specific to the new weekly record.
$ws.Cells.Item(181,4).Value  = 44642   # Fecha
$ws.Cells.Item(181,10).Value = 60      # Volumen
$ws.Cells.Item(181,11).Value = 550     # Precio minimo
$ws.Cells.Item(181,12).Value = 600     # Precio maximo
$ws.Cells.Item(181,13).Value = 575     # Precio promedio ponderado
$ws.Cells.Item(181,16).Value = 575     # Precio $/Kg
